# Apply weekly update: rotate data among rows 2,3,5,6,7 (row 4 stays fixed)
# Content moves: row2 -> row5, row5 -> row7, row7 -> row3, row3 -> row6, row6 -> row2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that contain data which changes across rows (D, K, L, M, N, O, P, Q, R, S, T)
$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

# Capture original values for the rows involved in the rotation
$orig = @{}
foreach ($r in 2,3,5,6,7) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

# Destination mapping: content that was in row X moves to row Y
$moveTo = @{ 2 = 5; 5 = 7; 7 = 3; 3 = 6; 6 = 2 }

foreach ($srcRow in $moveTo.Keys) {
    $destRow = $moveTo[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $orig[$srcRow][$c]
    }
}
